$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WorkSheet 1")

# Add the two trailing values to the existing row 9 (PriceChange, UpDown)
$ws.Range("X9").Value = -0.38999900000000309
$ws.Range("Y9").Value = "Down"

# Append a brand-new row 10 with a full record
$ws.Range("A10").Value = 42653.881261574075
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = "Buy"
$ws.Range("D10").Value = 30
$ws.Range("E10").Value = 12877
$ws.Range("F10").Value = 785
$ws.Range("G10").Value = 58
$ws.Range("H10").Value = 39
$ws.Range("I10").Value = 88
$ws.Range("J10").Value = 11
$ws.Range("K10").Value = 18011
$ws.Range("L10").Value = 130
$ws.Range("M10").Value = 88
$ws.Range("N10").Value = 48
$ws.Range("O10").Value = 6
$ws.Range("P10").Value = "Named"
$ws.Range("Q10").Value = 29.009773492518704
$ws.Range("R10").Value = 0.84
$ws.Range("S10").Value = -0.0136
$ws.Range("T10").Value = -0.03
$ws.Range("U10").Value = 14.53
$ws.Range("V10").Value = "N/A"
$ws.Range("W10").Value = -2

# Column widths refreshed by Excel's "best fit" recalculation after the new row was added
$ws.Columns.Item(1).ColumnWidth = 14.541666666666666
$ws.Columns.Item(2).ColumnWidth = 7.666666666666667
$ws.Columns.Item(3).ColumnWidth = 5.666666666666667
$ws.Columns.Item(4).ColumnWidth = 11.291666666666666
$ws.Columns.Item(5).ColumnWidth = 8.666666666666666
$ws.Columns.Item(6).ColumnWidth = 11.416666666666666
$ws.Columns.Item(7).ColumnWidth = 18.416666666666668
$ws.Columns.Item(8).ColumnWidth = 18.541666666666668
$ws.Columns.Item(9).ColumnWidth = 19.541666666666668
$ws.Columns.Item(10).ColumnWidth = 19.791666666666668
$ws.Columns.Item(11).ColumnWidth = 9.541666666666666
$ws.Columns.Item(12).ColumnWidth = 13.541666666666666
$ws.Columns.Item(13).ColumnWidth = 13.791666666666666

# Match the existing date / percentage formatting used by the rest of the table
$ws.Range("A10").NumberFormat = "m/d/yy h:mm"
$ws.Range("S10:T10").NumberFormat = "0.00%"
